# Apply the "Usuarios" -> "Clientes" template rework:
#  - rename the "Usuarios" sheet to "Clientes"
#  - update the form-name cell and a handful of field labels
#  - drop the "Clave" field row (the row below "Genero"/"Tipo cliente" data moves up)
#  - delete the now-unused "ZC_Login" sheet
#  - make "Clientes" the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- Rework "Usuarios" into "Clientes" -------------------------------------
$ws = $wb.Worksheets.Item("Usuarios")
$ws.Name = "Clientes"

$ws.Range("A2").Value = "Clientes"
$ws.Range("B4").Value = "Nombre"
$ws.Range("B6").Value = "Tipo cliente"
$ws.Range("K6").Value = "F=Natural,M=Juridica"
$ws.Range("B9").Value = "Fecha registro"

# Remove the "Clave" field row entirely (row 10); "Observaciones" shifts up.
$ws.Rows.Item(10).Delete()

# --- Remove the obsolete login-form sheet -----------------------------------
$wb.Worksheets.Item("ZC_Login").Delete()

# --- Make "Clientes" the active tab/selection -------------------------------
$ws.Activate()
$ws.Range("E16").Select()
